$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.250.50"

$ws.Range("E2").Value = "  +11.71%  "

$ws.Range("D3").Value = "1.684.36"

$ws.Range("E3").Value = "  +7.62%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.79"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  +9.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  +1.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3743"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = "  +1.96%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3463"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  +6.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.85"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "  +16.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.194"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  +6.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07330"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "  +4.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.55"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  +3.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.134"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  +6.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.803"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  +5.35%  "

$ws.Range("D16").Value = "1.685.60"

$ws.Range("E16").Value = "  +7.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001113"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  +4.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9989"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  +1.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06738"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  +9.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "82.40"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  +11.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.53"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  +3.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.131"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  +5.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.08"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "  +4.76%  "

$ws.Range("D24").Value = "24.179.53"

$ws.Range("E24").Value = "  +11.16%  "

$ws.Range("E25").Value = "  +2.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.694"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "  +13.50%  "

$ws.Range("E27").Value = "  -8.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.62"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = "  +3.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.71"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  +8.70%  "

$ws.Range("D30").Value = "1.869.87"

$ws.Range("E30").Value = "  +7.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.08"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = "  +6.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.490"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  +21.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.097"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  +0.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9941"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  +11.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.795"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  +12.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08490"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  +4.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.54"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  +9.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06501"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  +8.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.406"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = "  +6.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.951"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  +11.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02365"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  +9.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.280"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "  +3.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2143"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "  +8.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6230"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  +9.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9982"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  +2.26%  "

$ws.Range("B46").Value = "EnergySwap"

$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.36"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  +4.43%  "

$ws.Range("B47").Value = "PancakeSwap"

$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.812"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "  +5.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9987"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  +7.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.97"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  +3.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.042"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  +7.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07195"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = "  +7.09%  "
